$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for a99d8ed4-...md is now "Ready for handoff" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn detail sheet: new handoff timestamps + status update ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-02-18 04:17:17"
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-02-18 04:17:17"

# --- de-de detail sheet: new handoff timestamps + status update ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-02-18 04:17:30"
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-02-18 04:17:30"
